# Correct logging system configuration
# Appends one new data row to each of the four worksheets, mirroring the
# formatting of the existing rows directly above the new one.

$wb = $excel.ActiveWorkbook

function Add-LogRow($ws, $row, $aValue, $bValue, $cValue, $dValue, $eValue, $fValue, $gValue, $gIsText, $hValue, $iValue) {
    # Column A: date/time value, copy the number format used by the row above
    $ws.Cells.Item($row, 1).Value = [double]$aValue
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

    # Columns B-E: hex-ish text payloads, always stored as plain text
    $ws.Cells.Item($row, 2).Value = $bValue
    $ws.Cells.Item($row, 3).Value = $cValue
    $ws.Cells.Item($row, 4).Value = $dValue
    $ws.Cells.Item($row, 5).Value = $eValue

    # Column F: numeric length
    $ws.Cells.Item($row, 6).Value = $fValue

    # Column G: decoded ID - usually numeric, but occasionally too large for
    # a double to round-trip exactly, in which case it must stay text.
    if ($gIsText) {
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $gValue
        $ws.Cells.Item($row, 7).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 7).Value = [double]$gValue
    }

    # Columns H-I: numeric tallies
    $ws.Cells.Item($row, 8).Value = $hValue
    $ws.Cells.Item($row, 9).Value = $iValue
}

# --- Sheet "ROW50-FE-LIFTER" : add row 49 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-LogRow $ws1 49 "45748.687055" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x66" "0xe" 400 "5.68631262647114e+23" $false 358 14

# --- Sheet "ROW50-MID-LIFTER" : add row 51 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-LogRow $ws2 51 "45748.65979166667" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x6a" "0x19" 400 "568631262647113771663628" $true 362 25

# --- Sheet "ROW11-FE-LIFTER" : add row 49 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-LogRow $ws3 49 "45748.71592201389" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x66" "0x14" 400 "5.68631262647114e+23" $false 358 20

# --- Sheet "ROW11-MID-LIFTER" : add row 49 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-LogRow $ws4 49 "45748.85342299769" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x6e" "0x19" 400 "5.68631262647114e+23" $false 366 25
